$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the D_heel, D_h, D_yaw blocks (rows 7-9) and lift, drag, Frh blocks (rows 11-13) as completed
$ws.Range("G7").Value = "yes"
$ws.Range("G8").Value = "yes"
$ws.Range("G9").Value = "yes"
$ws.Range("G11").Value = "yes"
$ws.Range("G12").Value = "yes"
$ws.Range("G13").Value = "yes"

# Update the active cell selection
$ws.Range("I13").Select()
